$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit re-shuffles the observation rows 4-12 (row 7 is untouched) so that
# each row's identifying/species columns (A,B,D,E,F,G,H,Q,R) now hold the
# data that used to live in a different row, following one 8-cycle:
#   4<-5, 5<-10, 6<-9, 8<-6, 9<-12, 10<-8, 11<-4, 12<-11
# Row 10 also gains the empty K/L/M/N cells and the "hack" AC comment that
# row 8 used to carry (those stay with the observation, not the row number).

function Set-ObsRow($row, $A, $B, $D, $E, $F, $G, $H, $Q, $R) {
    $ws.Cells.Item($row, 1).Value = $A   # A - Id
    $ws.Cells.Item($row, 2).Value = $B   # B - Taxonsorteringsordning
    $ws.Cells.Item($row, 4).Value = $D   # D - Rodlistade
    $ws.Cells.Item($row, 5).Value = $E   # E - TaxonId
    $ws.Cells.Item($row, 6).Value = $F   # F - Artnamn
    $ws.Cells.Item($row, 7).Value = $G   # G - Vetenskapligt namn
    $ws.Cells.Item($row, 8).Value = $H   # H - Auktor
    $ws.Cells.Item($row, 17).Value = $Q  # Q - Ost
    $ws.Cells.Item($row, 18).Value = $R  # R - Nord
}

# New values for each destination row, taken from the source row noted above.
Set-ObsRow 4  111936768 90235 "LC" 3298   "Trådticka"  "Climacocystis borealis"     "(Fr.) Kotl. & Pouzar"                  490317 7088522
Set-ObsRow 5  111936779 77650 "NT" 6425   "Garnlav"    "Alectoria sarmentosa"       "(Ach.) Ach."                           490008 7088597
Set-ObsRow 6  111936781 89941 "LC" 4217   "Blodticka"  "Meruliopsis taxicola"       "(Pers.:Fr.) Bondartsev"                490315 7088552
Set-ObsRow 8  111936775 89567 "NT" 1204   "Gränsticka" "Phellopilus nigrolimitatus" "(Romell) Niemelä, T.Wagner & M.Fisch." 490380 7088379
Set-ObsRow 9  111936767 90235 "LC" 3298   "Trådticka"  "Climacocystis borealis"     "(Fr.) Kotl. & Pouzar"                  490377 7088412
Set-ObsRow 10 111936774 56446 "NT" 100049 "Spillkråka" "Dryocopus martius"          "(Linnaeus, 1758)"                      490378 7088551
Set-ObsRow 11 111936777 77650 "NT" 6425   "Garnlav"    "Alectoria sarmentosa"       "(Ach.) Ach."                           490056 7088709
Set-ObsRow 12 111936776 77650 "NT" 6425   "Garnlav"    "Alectoria sarmentosa"       "(Ach.) Ach."                           490398 7088445

# Move the empty Alder-Stadium/Kon/Aktivitet/Metod (K:N) cells off row 8 ...
$ws.Range("K8:N8").ClearContents()

# ... and onto row 10 (re-applying the default style keeps the cells present
# but empty, instead of removing them outright), alongside the "hack" public
# comment that travels with this particular observation.
$ws.Range("K10:N10").Style = "Normal"
$ws.Range("AC10").Value = "hack"

# The comment no longer belongs on row 8.
$ws.Range("AC8").ClearContents()
